# Fix 0.1mm misalignment of mounting holes:
# The BOM row for J4 (OLED module) had "N.M." placeholders in the LCSC (H)
# and MOUSER (I) columns; correct these to "-" to match the convention
# used by the other "Not Mounted" / no-part rows in the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BassFly_uHat.bom")

$ws.Range("H9").Value = "-"
$ws.Range("I9").Value = "-"

# Move the active selection to H15, matching the author's final cursor spot.
$ws.Range("H15").Select()
